$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.342.05"
$ws.Range("E2").Value = "  -2.26%  "
$ws.Range("D3").Value = "3.742.59"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.07%  "
$ws.Range("E6").Value = "  -4.01%  "
$ws.Range("D7").Value = "3.742.15"
$ws.Range("E7").Value = "  +2.14%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.531"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.10%  "
$ws.Range("E11").Value = "  -5.58%  "
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000256"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("D15").Value = "4.362.79"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "3.742.78"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").Value = "69.407.24"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "504.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.47%  "
$ws.Range("E23").Value = "  -2.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000137"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +23.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("E31").Value = "  +4.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.114"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.94%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  +4.38%  "
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.332"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.87%  "
$ws.Range("E39").Value = "  +2.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "44.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "420.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.14%  "
$ws.Range("D45").Value = "3.011.42"
$ws.Range("E45").Value = "  -3.79%  "
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("E47").Value = "  -2.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.71%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("E51").Value = "  -0.63%  "
